$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "SYMBOL_2000_SS_ENTRY_DATE"
$ws.Range("B4").Value = "SYMBOL_2000_SS_ENTRY_DATE"
$ws.Range("B5").Value = "SYMBOL_2000_SS_ENTRY_DATE"

$ws.Range("AI5").Value = 20190101

$ws.Range("Z3").Select()
